$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.3
$summary.Range("B4").Value = 0.29
$summary.Range("B6").Value = 67
$summary.Range("B8").Value = 35
$summary.Range("B9").Value = 31.34

# --- Strategy Status sheet (MarketMaking row, row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.3
$status.Range("D4").Value = 67
$status.Range("E4").Value = 0.29
$status.Range("F4").Value = 0.3
$status.Range("G4").Value = 31.34

# --- All Trades sheet: append new trade row 68 ---
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(68, 1).Value = 67
$allTrades.Cells.Item(68, 2).Value = "2026-02-17"
$allTrades.Cells.Item(68, 3).Value = "15:47:49"
$allTrades.Cells.Item(68, 4).Value = "MarketMaking"
$allTrades.Cells.Item(68, 5).Value = "DOWN"
$allTrades.Cells.Item(68, 6).Value = 0.25
$allTrades.Cells.Item(68, 7).Value = 0.24
$allTrades.Cells.Item(68, 8).Value = "CLOSED"
$allTrades.Cells.Item(68, 9).Value = -4
$allTrades.Cells.Item(68, 10).Value = -0.01
$allTrades.Cells.Item(68, 11).Value = 100.3
$allTrades.Cells.Item(68, 12).Value = 0
$allTrades.Cells.Item(68, 13).Value = 0
$allTrades.Cells.Item(68, 14).Value = 0.6
$allTrades.Cells.Item(68, 15).Value = "Normal spread capture: 19600 bps"
$allTrades.Cells.Item(68, 16).Value = "early_exit"
$allTrades.Cells.Item(68, 17).Value = 0.11

# --- MarketMaking sheet: append new trade row 68 (identical data) ---
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(68, 1).Value = 67
$mm.Cells.Item(68, 2).Value = "2026-02-17"
$mm.Cells.Item(68, 3).Value = "15:47:49"
$mm.Cells.Item(68, 4).Value = "MarketMaking"
$mm.Cells.Item(68, 5).Value = "DOWN"
$mm.Cells.Item(68, 6).Value = 0.25
$mm.Cells.Item(68, 7).Value = 0.24
$mm.Cells.Item(68, 8).Value = "CLOSED"
$mm.Cells.Item(68, 9).Value = -4
$mm.Cells.Item(68, 10).Value = -0.01
$mm.Cells.Item(68, 11).Value = 100.3
$mm.Cells.Item(68, 12).Value = 0
$mm.Cells.Item(68, 13).Value = 0
$mm.Cells.Item(68, 14).Value = 0.6
$mm.Cells.Item(68, 15).Value = "Normal spread capture: 19600 bps"
$mm.Cells.Item(68, 16).Value = "early_exit"
$mm.Cells.Item(68, 17).Value = 0.11
